# Client Setup Admin Suite - apply end user integrated test results
#
# - "Automation Supplier" test data is renamed to "Auto Supplier" (and its
#   short alias "aus" becomes "aut") across the ClientSupplierCreation,
#   ApplicationProvisioning and DashboardSetup sheets.
# - The Results columns for those rows are filled in with the outcome of the
#   integrated run (row 2 = SKIP, row 3 = PASS) and the "Test Cases" sheet's
#   Results column is marked PASS for every scenario.
# - DashboardSetup becomes the active/selected sheet, with B3 selected on it
#   (as well as on ClientSupplierCreation and ApplicationProvisioning).

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsClientSupplierCreation = $wb.Worksheets.Item("ClientSupplierCreation")
$wsApplicationProvisioning = $wb.Worksheets.Item("ApplicationProvisioning")
$wsDashboardSetup = $wb.Worksheets.Item("DashboardSetup")

# --- ClientSupplierCreation: clientAlias "aus" -> "aut" ---
$wsClientSupplierCreation.Range("D3").Value = "aut"

# --- Rename "Automation Supplier" -> "Auto Supplier" everywhere it is used ---
$wsClientSupplierCreation.Range("B3").Value = "Auto Supplier"
$wsApplicationProvisioning.Range("B3").Value = "Auto Supplier"
$wsDashboardSetup.Range("B3").Value = "Auto Supplier"

# --- Results: row 2 (header/blank data row) marked SKIP ---
$wsClientSupplierCreation.Range("K2").Value = "SKIP"
$wsApplicationProvisioning.Range("D2").Value = "SKIP"
$wsDashboardSetup.Range("D2").Value = "SKIP"

# --- Results: row 3 (actual test row) marked PASS ---
$wsClientSupplierCreation.Range("K3").Value = "PASS"
$wsApplicationProvisioning.Range("D3").Value = "PASS"
$wsDashboardSetup.Range("D3").Value = "PASS"
$wsTestCases.Range("D3").Value = "PASS"
$wsTestCases.Range("D4").Value = "PASS"
$wsTestCases.Range("D5").Value = "PASS"

# --- Selection / active sheet state ---
$wsClientSupplierCreation.Activate()
$wsClientSupplierCreation.Range("B3").Select()

$wsApplicationProvisioning.Activate()
$wsApplicationProvisioning.Range("B3").Select()

$wsDashboardSetup.Activate()
$wsDashboardSetup.Range("B3").Select()
